# Update NATMI LR-pair output sheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.152264333333333
$ws.Range("H2").Value = 3.456793
$ws.Range("I2").Value = 0.3930660006090215
$ws.Range("J2").Value = 0.3930660006090216
$ws.Range("M2").Value = 15.75563966666667
$ws.Range("N2").Value = 47.266919
$ws.Range("O2").Value = 0.3220556913988901
$ws.Range("P2").Value = 0.32205569139889
$ws.Range("Q2").Value = 18.15466163675189
$ws.Range("R2").Value = 163.391954730767
$ws.Range("S2").Value = 0.126589142591535
$ws.Range("T2").Value = 0.126589142591535

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.152264333333333
$ws.Range("H3").Value = 3.456793
$ws.Range("I3").Value = 0.3930660006090215
$ws.Range("J3").Value = 0.3930660006090216
$ws.Range("O3").Value = 0.5509544596378365
$ws.Range("P3").Value = 0.5509544596378364
$ws.Range("Q3").Value = 31.05795692831177
$ws.Range("R3").Value = 279.521612354806
$ws.Range("S3").Value = 0.216561465967549
$ws.Range("T3").Value = 0.216561465967549

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.152264333333333
$ws.Range("H4").Value = 3.456793
$ws.Range("I4").Value = 0.3930660006090215
$ws.Range("J4").Value = 0.3930660006090216
$ws.Range("O4").Value = 0.1269898489632735
$ws.Range("P4").Value = 0.1269898489632735
$ws.Range("Q4").Value = 7.158568535821889
$ws.Range("R4").Value = 64.42711682239701
$ws.Range("S4").Value = 0.04991539204993763
$ws.Range("T4").Value = 0.04991539204993762

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.4761983545501621
$ws.Range("J5").Value = 0.476198354550162
$ws.Range("M5").Value = 15.75563966666667
$ws.Range("N5").Value = 47.266919
$ws.Range("O5").Value = 0.3220556913988901
$ws.Range("P5").Value = 0.32205569139889
$ws.Range("Q5").Value = 21.99432152727834
$ws.Range("R5").Value = 197.948893745505
$ws.Range("S5").Value = 0.1533623903176662
$ws.Range("T5").Value = 0.1533623903176662

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.4761983545501621
$ws.Range("J6").Value = 0.476198354550162
$ws.Range("O6").Value = 0.5509544596378365
$ws.Range("P6").Value = 0.5509544596378364
$ws.Range("S6").Value = 0.2623636071116114
$ws.Range("T6").Value = 0.2623636071116113

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.4761983545501621
$ws.Range("J7").Value = 0.476198354550162
$ws.Range("O7").Value = 0.1269898489632735
$ws.Range("P7").Value = 0.1269898489632735
$ws.Range("S7").Value = 0.06047235712088446
$ws.Range("T7").Value = 0.06047235712088444

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.1307356448408163
$ws.Range("J8").Value = 0.1307356448408163
$ws.Range("M8").Value = 15.75563966666667
$ws.Range("N8").Value = 47.266919
$ws.Range("O8").Value = 0.3220556913988901
$ws.Range("P8").Value = 0.32205569139889
$ws.Range("Q8").Value = 6.038327894730444
$ws.Range("R8").Value = 54.344951052574
$ws.Range("S8").Value = 0.04210415848968884
$ws.Range("T8").Value = 0.04210415848968883

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.1307356448408163
$ws.Range("J9").Value = 0.1307356448408163
$ws.Range("O9").Value = 0.5509544596378365
$ws.Range("P9").Value = 0.5509544596378364
$ws.Range("S9").Value = 0.07202938655867608
$ws.Range("T9").Value = 0.07202938655867606

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.1307356448408163
$ws.Range("J10").Value = 0.1307356448408163
$ws.Range("O10").Value = 0.1269898489632735
$ws.Range("P10").Value = 0.1269898489632735
$ws.Range("S10").Value = 0.01660209979245144
$ws.Range("T10").Value = 0.01660209979245144
